$wb = $excel.ActiveWorkbook

# Worksheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value2 = 503
$ws.Range("I4").Value2 = 379.5
$ws.Range("J4").Value2 = 750
$ws.Range("K4").Value2 = 379.5
$ws.Range("L4").Value2 = 750
$ws.Range("M4").Value2 = -265.5
$ws.Range("N4").Value2 = -978
$ws.Range("H17").Value2 = 6546.6665
$ws.Range("J17").Value2 = 7109.4736
$ws.Range("L17").Value2 = 21328.4208
$ws.Range("N17").Value2 = -21664.4208
$ws.Range("H42").Value2 = 270.45456
$ws.Range("I42").Value2 = 202
$ws.Range("K42").Value2 = 606
$ws.Range("M42").Value2 = -376
$ws.Range("H137").Value2 = 2698.1643
$ws.Range("I137").Value2 = 1764.2354
$ws.Range("J137").Value2 = 2981.6785
$ws.Range("K137").Value2 = 5292.706200000001
$ws.Range("L137").Value2 = 8945.0355
$ws.Range("M137").Value2 = -2742.706200000001
$ws.Range("N137").Value2 = -14045.0355

# Worksheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 11026.25
$ws.Range("I2").Value2 = 1868.3334
$ws.Range("K2").Value2 = 1868.3334
$ws.Range("M2").Value2 = -1755.3334
$ws.Range("H44").Value2 = 20000
$ws.Range("J44").Value2 = 20000
$ws.Range("L44").Value2 = 20000
$ws.Range("N44").Value2 = -20976
$ws.Range("H61").Value2 = 9409
$ws.Range("I61").Value2 = 7673.615
$ws.Range("K61").Value2 = 7673.615
$ws.Range("M61").Value2 = -7461.615
$ws.Range("H74").Value2 = 16669776
$ws.Range("I74").Value2 = 41669356
$ws.Range("J74").Value2 = 3389.5
$ws.Range("K74").Value2 = 41669356
$ws.Range("L74").Value2 = 3389.5
$ws.Range("M74").Value2 = -41668482
$ws.Range("N74").Value2 = -5137.5
$ws.Range("H77").Value2 = 16669776
$ws.Range("I77").Value2 = 41669356
$ws.Range("J77").Value2 = 3389.5
$ws.Range("K77").Value2 = 208346780
$ws.Range("L77").Value2 = 16947.5
$ws.Range("M77").Value2 = -208342412
$ws.Range("N77").Value2 = -25683.5
$ws.Range("H80").Value2 = 105665.164
$ws.Range("I80").Value2 = 78000
$ws.Range("J80").Value2 = 133330.33
$ws.Range("K80").Value2 = 78000
$ws.Range("L80").Value2 = 133330.33
$ws.Range("M80").Value2 = -77002
$ws.Range("N80").Value2 = -135326.33
$ws.Range("H83").Value2 = 105665.164
$ws.Range("I83").Value2 = 78000
$ws.Range("J83").Value2 = 133330.33
$ws.Range("K83").Value2 = 234000
$ws.Range("L83").Value2 = 399990.99
$ws.Range("M83").Value2 = -229008
$ws.Range("N83").Value2 = -409974.99
$ws.Range("H116").Value2 = 11026.25
$ws.Range("I116").Value2 = 1868.3334
$ws.Range("K116").Value2 = 1868.3334
$ws.Range("M116").Value2 = 425.6666
$ws.Range("H132").Value2 = 2712.2263
$ws.Range("I132").Value2 = 2032.9778
$ws.Range("K132").Value2 = 6098.9334
$ws.Range("M132").Value2 = -3568.9334
$ws.Range("H136").Value2 = 9409
$ws.Range("I136").Value2 = 7673.615
$ws.Range("K136").Value2 = 23020.845
$ws.Range("M136").Value2 = -20470.845

# Worksheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 11026.25
$ws.Range("I3").Value2 = 1868.3334
$ws.Range("K3").Value2 = 1868.3334
$ws.Range("M3").Value2 = -1754.3334
$ws.Range("H20").Value2 = 3671
$ws.Range("I20").Value2 = 3329.5264
$ws.Range("J20").Value2 = 4211.6665
$ws.Range("K20").Value2 = 3329.5264
$ws.Range("L20").Value2 = 4211.6665
$ws.Range("M20").Value2 = -3082.5264
$ws.Range("N20").Value2 = -4705.6665
$ws.Range("H105").Value2 = 12207.032
$ws.Range("I105").Value2 = 11854.421
$ws.Range("K105").Value2 = 11854.421
$ws.Range("M105").Value2 = -10107.421

# Worksheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 10245.5
$ws.Range("I22").Value2 = 9986.5
$ws.Range("J22").Value2 = 10375
$ws.Range("K22").Value2 = 9986.5
$ws.Range("L22").Value2 = 10375
$ws.Range("M22").Value2 = -9636.5
$ws.Range("N22").Value2 = -11075
$ws.Range("H50").Value2 = 40000
$ws.Range("J50").Value2 = 40000
$ws.Range("L50").Value2 = 40000
$ws.Range("N50").Value2 = -41250
$ws.Range("H56").Value2 = 38750
$ws.Range("J56").Value2 = 45000
$ws.Range("L56").Value2 = 45000
$ws.Range("N56").Value2 = -46690
$ws.Range("H132").Value2 = 2840.6487
$ws.Range("I132").Value2 = 2052.28
$ws.Range("K132").Value2 = 6156.84
$ws.Range("M132").Value2 = -3626.84
$ws.Range("H134").Value2 = 2462.6177
$ws.Range("I134").Value2 = 1688.0435
$ws.Range("K134").Value2 = 5064.1305
$ws.Range("M134").Value2 = -2529.1305

# Worksheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value2 = 1309.8
$ws.Range("I17").Value2 = 633
$ws.Range("J17").Value2 = 2325
$ws.Range("K17").Value2 = 1899
$ws.Range("L17").Value2 = 6975
$ws.Range("M17").Value2 = -1730
$ws.Range("N17").Value2 = -7313
$ws.Range("H134").Value2 = 6642.8335
$ws.Range("I134").Value2 = 5750.684
$ws.Range("J134").Value2 = 10033
$ws.Range("K134").Value2 = 17252.052
$ws.Range("L134").Value2 = 30099
$ws.Range("M134").Value2 = -12182.052
$ws.Range("N134").Value2 = -40239

# Worksheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 860.61536
$ws.Range("I2").Value2 = 198.83333
$ws.Range("K2").Value2 = 198.83333
$ws.Range("M2").Value2 = -85.83332999999999
$ws.Range("H107").Value2 = 875.7273
$ws.Range("J107").Value2 = 760
$ws.Range("L107").Value2 = 760
$ws.Range("N107").Value2 = -4600
$ws.Range("H122").Value2 = 7652.517
$ws.Range("I122").Value2 = 6276.533
$ws.Range("K122").Value2 = 18829.599
$ws.Range("M122").Value2 = -16379.599

# Worksheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 8911.714
$ws.Range("I7").Value2 = 5669.4546
$ws.Range("K7").Value2 = 5669.4546
$ws.Range("M7").Value2 = -5557.4546
$ws.Range("H46").Value2 = 6753.75
$ws.Range("I46").Value2 = 5076.6
$ws.Range("J46").Value2 = 7516.091
$ws.Range("K46").Value2 = 5076.6
$ws.Range("L46").Value2 = 7516.091
$ws.Range("M46").Value2 = -4888.6
$ws.Range("N46").Value2 = -7892.091
$ws.Range("H93").Value2 = 1870.8125
$ws.Range("I93").Value2 = 1795.6
$ws.Range("K93").Value2 = 1795.6
$ws.Range("M93").Value2 = -547.5999999999999
$ws.Range("H122").Value2 = 374027
$ws.Range("I122").Value2 = 579700
$ws.Range("J122").Value2 = 14099.25
$ws.Range("K122").Value2 = 1739100
$ws.Range("L122").Value2 = 42297.75
$ws.Range("M122").Value2 = -1736650
$ws.Range("N122").Value2 = -47197.75
$ws.Range("H126").Value2 = 8911.714
$ws.Range("I126").Value2 = 5669.4546
$ws.Range("K126").Value2 = 17008.3638
$ws.Range("M126").Value2 = -14538.3638
$ws.Range("H135").Value2 = 89199.2
$ws.Range("J135").Value2 = 89199.2
$ws.Range("L135").Value2 = 89199.2
$ws.Range("N135").Value2 = -99339.2

# Worksheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 1773.5714
$ws.Range("I122").Value2 = 1202.8948
$ws.Range("K122").Value2 = 3608.6844
$ws.Range("M122").Value2 = -1158.6844
